# Add new jewelry products (Style 2172-2179) to Sheet1, rows 3-10,
# mirroring the format of the existing row 2 (style/format carried down).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New product rows: Style, Description, Color, Retail, SzRange (E).
# Column F stays blank and column G reuses the fixed care-instruction text,
# exactly like row 2.
$rows = @(
    @("2172", "Devotion Necklace",   "Pearl",       119),
    @("2173", "Palm Beach Earrings", "Seashell",     36),
    @("2174", "Palm Beach Necklace", "Seashell",     89),
    @("2175", "Palm Beach Bracelet", "Seashell",     79),
    @("2176", "La Bohème Earrings",  "Terra Cotta",  49),
    @("2177", "La Bohème Necklace",  "Terra Cotta",  79),
    @("2178", "Curio Earrings",      "Jasper",       49),
    @("2179", "Curio Necklace",      "Jasper",      129)
)

# Duplicate row 2 (format + values) into rows 3..10, pushing nothing else
# down since the sheet currently ends at row 2.
for ($i = 0; $i -lt $rows.Count; $i++) {
    $ws.Rows("2").Copy()
    $ws.Rows("3").Insert()
}

# Fill in the per-row data. Column A ("Style") holds numeric-looking codes
# (e.g. "2172") that must stay TEXT (matching the existing convention used
# for every other Style value in this sheet/workbook) rather than being
# auto-coerced to a number. Round-tripping the digits through a TEXT()
# formula and pasting the result as a value keeps the cell's text type
# (and its existing style) without minting a new number format/style.
$r = 3
foreach ($row in $rows) {
    $style = $row[0]
    $desc = $row[1]
    $color = $row[2]
    $retail = $row[3]

    $ws.Range("ZZ1").Formula = "=TEXT(" + $style + ",""0"")"
    $ws.Range("ZZ1").Copy()
    $ws.Range("A$r").PasteSpecial(-4163)  # xlPasteValues
    $ws.Range("ZZ1").Clear()

    $ws.Range("B$r").Value = $desc
    $ws.Range("C$r").Value = $color
    $ws.Range("D$r").Value = $retail
    # Column E (SzRange) and G (CareInstr) already carried "One Size" / the
    # care text down from row 2 via the row copy above - nothing to do.

    $r = $r + 1
}

$ws.Range("A2").Select()
$ws.Range("A2:XFD10").Select()
